# Updated symbol list on Mon Jan  2 03:52:42 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) columns on the active
# sheet hold plain text (e.g. "242.97", "-0.95%") rather than real numbers
# or percentages, so each value is written back as an explicit Text string
# (NumberFormat "@") to keep Excel from re-interpreting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "243.04"
    "E2"  = "-0.92%"
    "D3"  = "29.93"
    "E3"  = "12.82%"
    "D4"  = "5.125"
    "E4"  = "0.01%"
    "E5"  = "1.17%"
    "D6"  = "6.497"
    "E6"  = "0.33%"
    "D7"  = "0.8260"
    "E7"  = "1.07%"
    "D8"  = "0.8622"
    "E8"  = "3.28%"
    "D9"  = "0.1331"
    "E9"  = "-0.10%"
    "D10" = "0.06877"
    "E10" = "-1.67%"
    "D11" = "0.02856"
    "E11" = "-1.02%"
    "D12" = "0.09389"
    "D13" = "0.001509"
    "E13" = "-0.51%"
    "D14" = "0.04142"
    "E14" = "-9.70%"
    "D15" = "0.0006026"
    "E15" = "0.78%"
    "D16" = "0.006215"
    "E16" = "1.66%"
    "E17" = "-3.49%"
    "D18" = "3.012"
    "E18" = "-0.86%"
    "D19" = "2.215"
    "E19" = "1.49%"
    "E20" = "1.20%"
    "D21" = "0.03241"
    "E21" = "3.27%"
    "E22" = "-0.37%"
    "E23" = "-3.61%"
    "E24" = "-0.06%"
    "E25" = "-3.12%"
    "D26" = "0.004444"
    "E26" = "-1.58%"
    "E27" = "22.80%"
    "D40" = "0.03701"
    "E40" = "1.65%"
    "D41" = "0.005753"
    "E41" = "-7.05%"
    "E42" = "0.30%"
    "D43" = "0.002309"
    "E43" = "-3.76%"
    "D44" = "0.009658"
    "E44" = "8.86%"
    "D45" = "0.00005098"
    "E45" = "-4.43%"
    "D46" = "0.00000000750"
    "E46" = "-0.05%"
    "D48" = "0.002651"
    "E48" = "3.86%"
    "E49" = "-0.05%"
    "E50" = "-0.05%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
